$wb = $excel.ActiveWorkbook

# Sheet4: bump the base test identities from generation 50 to generation 52.
# All downstream sheets (Sheet2, Sheet3, Sheet5..Sheet8) reference these
# cells (directly or transitively) via formulas, so they will recalculate.
$sheet4 = $wb.Worksheets.Item("Sheet4")
$sheet4.Range("A2").Value = "tavalinetont52"
$sheet4.Range("C2").Value = "puhtaloom52"
$sheet4.Range("E2").Value = "filmweird52"

# Sheet5: bump the mailinator test accounts from generation 23 to generation 25.
$sheet5 = $wb.Worksheets.Item("Sheet5")
$sheet5.Range("I2").Value = "selentest25@mailinator.com"
$sheet5.Range("I3").Value = "testimeauto25@mailinator.com"
$sheet5.Range("I4").Value = "vahekonto25@mailinator.com"

# Sheet9: update the browser upload size test data (row 2 = firefox, row 3 = chrome).
$sheet9 = $wb.Worksheets.Item("Sheet9")
$sheet9.Range("C2").Value = "2309"
$sheet9.Range("D2").Value = "6956"
$sheet9.Range("E2").Value = "8025"
$sheet9.Range("G2").Value = "1356"
$sheet9.Range("H2").Value = "1309"
$sheet9.Range("I2").Value = "1254"

$sheet9.Range("C3").Value = "1782"
$sheet9.Range("D3").Value = "7262"
$sheet9.Range("E3").Value = "7742"
$sheet9.Range("G3").Value = "1027"
$sheet9.Range("H3").Value = "1090"
$sheet9.Range("I3").Value = "1067"
